$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the merged-cell remnants in D2:E2 (part of the C2:E2 merge) so the
# cells carry no leftover style/content, matching the cleaned-up header row.
$ws.Range("D2:E2").ClearFormats()
$ws.Range("D2:E2").ClearContents()

# New sequence of ESTIMATE DATE values for rows 4-93 (2025-08-01 .. 2025-11-14).
$newDates = @("2025-08-01","2025-08-02","2025-08-03","2025-08-04","2025-08-05","2025-08-06","2025-08-07","2025-08-08","2025-08-09","2025-08-10","2025-08-11","2025-08-12","2025-08-13","2025-08-14","2025-08-15","2025-08-16","2025-08-17","2025-08-18","2025-08-19","2025-08-20","2025-08-21","2025-08-22","2025-08-23","2025-08-24","2025-08-25","2025-08-26","2025-08-27","2025-08-28","2025-08-29","2025-08-30","2025-09-08","2025-09-09","2025-09-10","2025-09-11","2025-09-12","2025-09-13","2025-09-14","2025-09-15","2025-09-16","2025-09-17","2025-09-18","2025-09-19","2025-09-20","2025-09-21","2025-09-22","2025-09-23","2025-09-24","2025-09-25","2025-09-26","2025-09-27","2025-09-28","2025-09-29","2025-09-30","2025-10-01","2025-10-02","2025-10-03","2025-10-04","2025-10-05","2025-10-06","2025-10-07","2025-10-16","2025-10-17","2025-10-18","2025-10-19","2025-10-20","2025-10-21","2025-10-22","2025-10-23","2025-10-24","2025-10-25","2025-10-26","2025-10-27","2025-10-28","2025-10-29","2025-10-30","2025-10-31","2025-11-01","2025-11-02","2025-11-03","2025-11-04","2025-11-05","2025-11-06","2025-11-07","2025-11-08","2025-11-09","2025-11-10","2025-11-11","2025-11-12","2025-11-13","2025-11-14")

for ($i = 0; $i -lt $newDates.Count; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 8).Value = $newDates[$i]
}

# Dates should display without the time portion now (yyyy-mm-dd).
$ws.Range("H4:H93").NumberFormat = "yyyy-mm-dd"
